$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between row 4 and row 5 and need to be swapped.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "S", "Z", "AB")

foreach ($col in $cols) {
    $cell4 = $ws.Range($col + "4")
    $cell5 = $ws.Range($col + "5")

    $v4 = $cell4.Value2
    $v5 = $cell5.Value2

    $cell4.Value = $v5
    $cell5.Value = $v4
}
